# Fruta / hortaliza, semanal
# Insert two new weekly price rows (759-760) into the Coliflor dataset,
# shifting all the existing rows below down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 759.
$ws.Range("A759:A760").EntireRow.Insert()

# --- New row 759 -----------------------------------------------------
$ws.Cells.Item(759, 1).Value = 6
$ws.Cells.Item(759, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(759, 3).Value = "Metropolitana"
$ws.Cells.Item(759, 4).Value = 44753
$ws.Cells.Item(759, 5).Value = 13
$ws.Cells.Item(759, 6).Value = 100112008
$ws.Cells.Item(759, 7).Value = "Coliflor"
$ws.Cells.Item(759, 8).Value = "Sin especificar"
$ws.Cells.Item(759, 9).Value = "Primera"
$ws.Cells.Item(759, 10).Value = 6100
$ws.Cells.Item(759, 11).Value = 1000
$ws.Cells.Item(759, 12).Value = 1200
$ws.Cells.Item(759, 13).Value = 1111
$ws.Cells.Item(759, 14).Value = "`$/unidad"
$ws.Cells.Item(759, 15).Value = "Región Metropolitana"
$ws.Cells.Item(759, 16).Value = 1111
$ws.Cells.Item(759, 17).Value = 1
$ws.Cells.Item(759, 18).Value = "Hortaliza"

# --- New row 760 -----------------------------------------------------
$ws.Cells.Item(760, 1).Value = 6
$ws.Cells.Item(760, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(760, 3).Value = "Metropolitana"
$ws.Cells.Item(760, 4).Value = 44753
$ws.Cells.Item(760, 5).Value = 13
$ws.Cells.Item(760, 6).Value = 100112008
$ws.Cells.Item(760, 7).Value = "Coliflor"
$ws.Cells.Item(760, 8).Value = "Sin especificar"
$ws.Cells.Item(760, 9).Value = "Segunda"
$ws.Cells.Item(760, 10).Value = 2200
$ws.Cells.Item(760, 11).Value = 800
$ws.Cells.Item(760, 12).Value = 800
$ws.Cells.Item(760, 13).Value = 800
$ws.Cells.Item(760, 14).Value = "`$/unidad"
$ws.Cells.Item(760, 15).Value = "Región Metropolitana"
$ws.Cells.Item(760, 16).Value = 800
$ws.Cells.Item(760, 17).Value = 1
$ws.Cells.Item(760, 18).Value = "Hortaliza"
